$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 365
$ws.Range("I2").Value = 153.33333
$ws.Range("K2").Value = 153.33333
$ws.Range("M2").Value = -40.33332999999999
$ws.Range("H6").Value = 494
$ws.Range("I6").Value = 489
$ws.Range("J6").Value = 499
$ws.Range("K6").Value = 1467
$ws.Range("L6").Value = 1497
$ws.Range("M6").Value = -1355
$ws.Range("N6").Value = -1721
$ws.Range("H43").Value = 3856903
$ws.Range("I43").Value = 3856903
$ws.Range("K43").Value = 3856903
$ws.Range("M43").Value = -3856834
$ws.Range("H70").Value = 3132.4736
$ws.Range("J70").Value = 3200.8333
$ws.Range("L70").Value = 9602.499899999999
$ws.Range("N70").Value = -10142.4999
$ws.Range("H73").Value = 3132.4736
$ws.Range("J73").Value = 3200.8333
$ws.Range("L73").Value = 9602.499899999999
$ws.Range("N73").Value = -11474.4999
$ws.Range("H88").Value = 9097132
$ws.Range("J88").Value = 6900.3125
$ws.Range("L88").Value = 6900.3125
$ws.Range("N88").Value = -7712.3125
$ws.Range("H91").Value = 9097132
$ws.Range("J91").Value = 6900.3125
$ws.Range("L91").Value = 6900.3125
$ws.Range("N91").Value = -9708.3125
$ws.Range("H106").Value = 83338216
$ws.Range("I106").Value = 333333340
$ws.Range("J106").Value = 6502
$ws.Range("K106").Value = 333333340
$ws.Range("L106").Value = 6502
$ws.Range("M106").Value = -333332709
$ws.Range("N106").Value = -7764
$ws.Range("H107").Value = 628.8889
$ws.Range("I107").Value = 662
$ws.Range("K107").Value = 662
$ws.Range("M107").Value = 1258
$ws.Range("H125").Value = 9054.111000000001
$ws.Range("I125").Value = 4299.8
$ws.Range("J125").Value = 14997
$ws.Range("K125").Value = 38698.2
$ws.Range("L125").Value = 134973
$ws.Range("M125").Value = -36238.2
$ws.Range("N125").Value = -139893
$ws.Range("H132").Value = 12565.712
$ws.Range("I132").Value = 7621.3335
$ws.Range("K132").Value = 22864.0005
$ws.Range("M132").Value = -20334.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 646.6
$ws.Range("I5").Value = 646.6
$ws.Range("K5").Value = 646.6
$ws.Range("M5").Value = -534.6
$ws.Range("H32").Value = 4183.4126
$ws.Range("I32").Value = 2344.0408
$ws.Range("J32").Value = 10621.214
$ws.Range("K32").Value = 2344.0408
$ws.Range("L32").Value = 10621.214
$ws.Range("M32").Value = -2057.0408
$ws.Range("N32").Value = -11195.214
$ws.Range("H40").Value = 29000
$ws.Range("J40").Value = 29000
$ws.Range("L40").Value = 29000
$ws.Range("N40").Value = -29352
$ws.Range("H45").Value = 1679.6364
$ws.Range("I45").Value = 1469.6666
$ws.Range("K45").Value = 1469.6666
$ws.Range("M45").Value = -1092.6666
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30520
$ws.Range("H74").Value = 27779264
$ws.Range("J74").Value = 1554
$ws.Range("L74").Value = 1554
$ws.Range("N74").Value = -3302
$ws.Range("H77").Value = 27779264
$ws.Range("J77").Value = 1554
$ws.Range("L77").Value = 7770
$ws.Range("N77").Value = -16506
$ws.Range("H110").Value = 662121.9399999999
$ws.Range("I110").Value = 853807.3
$ws.Range("J110").Value = 4915.143
$ws.Range("K110").Value = 853807.3
$ws.Range("L110").Value = 4915.143
$ws.Range("M110").Value = -851762.3
$ws.Range("N110").Value = -9005.143
$ws.Range("H132").Value = 21828.531
$ws.Range("I132").Value = 24777.041
$ws.Range("K132").Value = 74331.12300000001
$ws.Range("M132").Value = -71801.12300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 646.6
$ws.Range("I4").Value = 646.6
$ws.Range("K4").Value = 646.6
$ws.Range("M4").Value = -531.6
$ws.Range("H22").Value = 300.85715
$ws.Range("I22").Value = 334.33334
$ws.Range("K22").Value = 334.33334
$ws.Range("M22").Value = -161.33334
$ws.Range("H64").Value = 8334132.5
$ws.Range("I64").Value = 10417366
$ws.Range("J64").Value = 1200
$ws.Range("K64").Value = 10417366
$ws.Range("L64").Value = 1200
$ws.Range("M64").Value = -10417141
$ws.Range("N64").Value = -1650
$ws.Range("H67").Value = 8334132.5
$ws.Range("I67").Value = 10417366
$ws.Range("J67").Value = 1200
$ws.Range("K67").Value = 10417366
$ws.Range("L67").Value = 1200
$ws.Range("M67").Value = -10416586
$ws.Range("N67").Value = -2760
$ws.Range("H86").Value = 40002356
$ws.Range("I86").Value = 2293.2
$ws.Range("J86").Value = 100002450
$ws.Range("K86").Value = 2293.2
$ws.Range("L86").Value = 100002450
$ws.Range("M86").Value = -1170.2
$ws.Range("N86").Value = -100004696
$ws.Range("H89").Value = 40002356
$ws.Range("I89").Value = 2293.2
$ws.Range("J89").Value = 100002450
$ws.Range("K89").Value = 11466
$ws.Range("L89").Value = 500012250
$ws.Range("M89").Value = -5850
$ws.Range("N89").Value = -500023482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1451
$ws.Range("I16").Value = 1351.7778
$ws.Range("J16").Value = 1897.5
$ws.Range("K16").Value = 1351.7778
$ws.Range("L16").Value = 1897.5
$ws.Range("M16").Value = -1064.7778
$ws.Range("N16").Value = -2471.5
$ws.Range("H58").Value = 1001704.3
$ws.Range("I58").Value = 1430332.8
$ws.Range("J58").Value = 1571.3334
$ws.Range("K58").Value = 1430332.8
$ws.Range("L58").Value = 1571.3334
$ws.Range("M58").Value = -1430129.8
$ws.Range("N58").Value = -1977.3334
$ws.Range("H113").Value = 1451
$ws.Range("I113").Value = 1351.7778
$ws.Range("J113").Value = 1897.5
$ws.Range("K113").Value = 1351.7778
$ws.Range("L113").Value = 1897.5
$ws.Range("M113").Value = 818.2221999999999
$ws.Range("N113").Value = -6237.5
$ws.Range("H122").Value = 2961570.5
$ws.Range("J122").Value = 2734.25
$ws.Range("L122").Value = 8202.75
$ws.Range("N122").Value = -13102.75
$ws.Range("H136").Value = 1001704.3
$ws.Range("I136").Value = 1430332.8
$ws.Range("J136").Value = 1571.3334
$ws.Range("K136").Value = 4290998.4
$ws.Range("L136").Value = 4714.0002
$ws.Range("M136").Value = -4288448.4
$ws.Range("N136").Value = -9814.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1310.8
$ws.Range("J5").Value = 3002.5
$ws.Range("L5").Value = 9007.5
$ws.Range("N5").Value = -9231.5
$ws.Range("H122").Value = 482.7857
$ws.Range("I122").Value = 369.33334
$ws.Range("J122").Value = 513.7273
$ws.Range("K122").Value = 3324.00006
$ws.Range("L122").Value = 4623.545700000001
$ws.Range("M122").Value = -874.0000600000003
$ws.Range("N122").Value = -9523.545700000001
$ws.Range("H135").Value = 1310.8
$ws.Range("J135").Value = 3002.5
$ws.Range("L135").Value = 27022.5
$ws.Range("N135").Value = -32092.5
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 386.27274
$ws.Range("I97").Value = 399.93332
$ws.Range("J97").Value = 357
$ws.Range("K97").Value = 399.93332
$ws.Range("L97").Value = 357
$ws.Range("M97").Value = 96.06668000000002
$ws.Range("N97").Value = -1349
$ws.Range("H113").Value = 905024.4399999999
$ws.Range("J113").Value = 7222.222
$ws.Range("L113").Value = 7222.222
$ws.Range("N113").Value = -11562.222
$ws.Range("H126").Value = 3866.348
$ws.Range("I126").Value = 2053.3333
$ws.Range("J126").Value = 7265.75
$ws.Range("K126").Value = 6159.999899999999
$ws.Range("L126").Value = 21797.25
$ws.Range("M126").Value = -3689.999899999999
$ws.Range("N126").Value = -26737.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1221.45
$ws.Range("I22").Value = 1069.8572
$ws.Range("K22").Value = 1069.8572
$ws.Range("M22").Value = -774.8571999999999
$ws.Range("H27").Value = 1221.45
$ws.Range("I27").Value = 1069.8572
$ws.Range("K27").Value = 1069.8572
$ws.Range("M27").Value = -962.8571999999999
$ws.Range("H40").Value = 7139.143
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864
$ws.Range("H122").Value = 3681.4285
$ws.Range("I122").Value = 3339.6365
$ws.Range("K122").Value = 10018.9095
$ws.Range("M122").Value = -7568.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null
$ws.Range("H96").Value = 5396.857
$ws.Range("I96").Value = 4451.25
$ws.Range("J96").Value = 6657.6665
$ws.Range("K96").Value = 4451.25
$ws.Range("L96").Value = 6657.6665
$ws.Range("M96").Value = -3078.25
$ws.Range("N96").Value = -9403.666499999999
$ws.Range("H122").Value = 2460.6843
$ws.Range("I122").Value = 2303.4893
$ws.Range("K122").Value = 6910.467900000001
$ws.Range("M122").Value = -4460.467900000001
$ws.Range("H126").Value = 1303.5
$ws.Range("I126").Value = 1297.125
$ws.Range("J126").Value = 1312
$ws.Range("K126").Value = 3891.375
$ws.Range("L126").Value = 3936
$ws.Range("M126").Value = -1421.375
$ws.Range("N126").Value = -8876

Write-Output "applied edits"